# Maddelavedu_LabExam03Grading.xlsx - grading update
# Prakash 33 - 59 (upto Maddelavedu_LabExam03Grading.xlsx)
#
# Grade two previously-ungraded rubric rows (14 and 16 -> worksheet rows 22 and 24)
# with full points (10), which cascades through the existing SUM() formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Award points for "addProduct() method" (row 22) and "whoPurchasedProduct() method" (row 24)
$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 10

# Recalculate so the dependent totals (E26, E38) pick up the new values
$excel.Calculate()

# Restore the view: scrolled so row 8 is at the top, with E24 as the active selection
$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E24").Select()
